$wb = $excel.ActiveWorkbook

# Insert a new worksheet "total_concentrations" right before
# "equilibrium_concentrations" (i.e. it becomes sheet #4, pushing
# equilibrium_concentrations / L_fractions / percent_error / component_names
# each one position later).
$target = $wb.Worksheets.Item("equilibrium_concentrations")
$ws = $wb.Worksheets.Add($target)
$ws.Name = "total_concentrations"

# Header row
$ws.Range("A1").Value = "H"
$ws.Range("B1").Value = "L"
$ws.Range("C1").Value = "M"

# Data rows
$ws.Range("A2").Value = 0.0282807977164644
$ws.Range("B2").Value = 0.02
$ws.Range("C2").Value = 0.0103840613673434

$ws.Range("A3").Value = 0.0110324008139916
$ws.Range("B3").Value = 0.02
$ws.Range("C3").Value = 0.0123740633160388

$ws.Range("A4").Value = 0.00114136198682271
$ws.Range("B4").Value = 0.01
$ws.Range("C4").Value = 0.00623870084795149
